# Apply saved-results updates to Sheet1 per commit "saved results for responses".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("S2").Value  = 418
$ws.Range("U3").Value  = 401
$ws.Range("V3").Value  = 421
$ws.Range("W3").Value  = 428
$ws.Range("X3").Value  = 458
$ws.Range("T6").Value  = 437
$ws.Range("U6").Value  = 429
$ws.Range("V6").Value  = 417
$ws.Range("T7").Value  = 436
$ws.Range("R10").Value = 416
$ws.Range("W10").Value = 417
$ws.Range("O13").Value = 417
$ws.Range("U17").Value = 428
